$d = $word.ActiveDocument

# Anchor: the pre-existing trailing empty paragraph (last paragraph in the doc).
$anchorPara = $d.Paragraphs($d.Paragraphs.Count)
$anchor = $d.Range($anchorPara.Range.Start, $anchorPara.Range.Start)

# Insert four new paragraphs before the anchor, in one shot, using placeholder
# markers we will replace with formatted runs / a hyperlink afterwards.
# Layout (after insertion, before anchor's own empty paragraph):
#   (empty)
#   Creator Assets. (2019, August 16). <<ITALIC>>. YouTube. <<URL>>
#   (empty)
#   <<ZWNJ>>
$anchor.InsertBefore("`r" + "Creator Assets. (2019, August 16). ITALICPLACEHOLDER. YouTube. URLPLACEHOLDER`r`rZWNJPLACEHOLDER`r")

# --- Italicize the title run ---
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute("ITALICPLACEHOLDER", $true, $true, $false, $false, $false, `
                   $true, 1, $false, "Click Sound Effects (Copyright Free)", 2)

$rng2 = $d.Content
$rng2.Find.ClearFormatting()
$rng2.Find.Execute("Click Sound Effects (Copyright Free)", $true, $true, $false, $false, $false, `
                    $true, 1, $false, "", 0)
$rng2.Font.Italic = $true
$rng2.Font.ItalicBi = $true

# --- Append ". YouTube. " literal already present; now turn the URL placeholder into a hyperlink ---
$urlRng = $d.Content
$urlRng.Find.ClearFormatting()
$urlRng.Find.Execute("URLPLACEHOLDER", $true, $true, $false, $false, $false, `
                      $true, 1, $false, "", 0)
$d.Hyperlinks.Add($urlRng, "https://www.youtube.com/watch?v=q8ZLBOFQ2g0") | Out-Null

# --- Replace the ZWNJ placeholder with an actual zero-width non-joiner char ---
$zwnjRng = $d.Content
$zwnjRng.Find.ClearFormatting()
$zwnjRng.Find.Execute("ZWNJPLACEHOLDER", $true, $true, $false, $false, $false, `
                       $true, 1, $false, [string][char]0x200C, 2)

Write-Output ("Paragraphs.Count=" + $d.Paragraphs.Count)
